$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = -0.06847097740271614; C = 0.6620077945808225;  D = 0.8576394255630952; E = 0.9260882385405266; F = 0.9347485423480986; G = 42 }
    3  = @{ B = 0.1628294350496699;   C = 0.6137740522511609;  D = 0.7803318254486972; E = 0.8833639258248535; F = 0.8790130036258593; G = 41 }
    4  = @{ B = -0.06854388644398214; C = 0.6834144119998282;  D = 0.8849781572300088; E = 0.9407327767384364; F = 0.9501848094793431; G = 40 }
    5  = @{ B = 0.1479423124587954;   C = 0.6621921511100777;  D = 0.852793974404251;  E = 0.9234684479744021; F = 0.9234571032261442; G = 39 }
    6  = @{ B = -0.08040621195020153; C = 0.6834744100144207;  D = 0.8476544914969599; E = 0.920681536415801;  F = 0.9294752172076342; G = 38 }
    7  = @{ B = 0.1165699234898814;   C = 0.6791032791661118;  D = 0.8276536103283203; E = 0.9097546978874692; F = 0.9147010340487955; G = 37 }
    8  = @{ B = -0.1318557986420442;  C = 0.6256777560484761;  D = 0.7095626088736116; E = 0.8423553934495888; F = 0.843773152665554;  G = 36 }
    9  = @{ B = 0.04870157088247548;  C = 0.6198098084369814;  D = 0.7254827416971069; E = 0.8517527468092527; F = 0.8627739369610303; G = 35 }
    10 = @{ B = -0.06798252598835164; C = 0.5847756573626377;  D = 0.6476270503415246; E = 0.8047527883403229; F = 0.8139351567699135; G = 34 }
    11 = @{ B = 0.01270493608304687;  C = 0.6241980251273542;  D = 0.7283763178587769; E = 0.853449657483543;  F = 0.8665861862271473; G = 33 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
